$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting (style) from BA1 (date header style) to the new BB1 header cell
$ws.Range("BA1").Copy()
$ws.Range("BB1").PasteSpecial(-4122)
$ws.Range("BB1").Value = 45986

$ws.Range("BB3").Value = -0.6303298271467694
$ws.Range("BB4").Value = 0.3824851463237522
$ws.Range("BB5").Value = 0.4942252260480062
$ws.Range("BB6").Value = -0.8283953561833202
$ws.Range("BB7").Value = -0.1541135436267549
$ws.Range("BB8").Value = 0.2992729818363626
$ws.Range("BB9").Value = -0.6011708148489947
$ws.Range("BB10").Value = -0.2513629445286991
$ws.Range("BB11").Value = 0.2729872858366011
$ws.Range("BB12").Value = 0.1243096661369014
$ws.Range("BB13").Value = -0.4025146932836732
$ws.Range("BB14").Value = -0.9275935716973494
$ws.Range("BB15").Value = 0.2398177392026746
$ws.Range("BB16").Value = 0.539177729005802
$ws.Range("BB17").Value = 0.6992203852886458
$ws.Range("BB18").Value = -0.850803046382087
$ws.Range("BB19").Value = 1.049317648994741
$ws.Range("BB20").Value = 0.3243937446859801
$ws.Range("BB21").Value = 0.1600960256025408
